$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("TN330", "Natalie's - Honey Tangerine", "1", "14.57", "14.57"),
    @("TN374", "Natalie's - Lemonade", "1", "9.30", "9.30"),
    @("AH252", "Natalie's - Orange Juice", "2", "24.50", "49.00"),
    @("TN454", "Natalie's - Orange Mango", "1", "13.38", "13.38"),
    @("TN362", "Natalie's - Orange Pineapple", "2", "13.38", "26.76"),
    @("TN380", "Natalie's - Strawberry Lemonade", "1", "10.15", "10.15")
)

$row = 20
foreach ($item in $data) {
    # Columns C, D, E hold digit-only strings ("1", "14.57", ...). Plain
    # .Value assignment would let them be auto-coerced to numbers (and in
    # the process mangle values like 9.30/13.38 into binary-float noise),
    # so force a Text number format first to keep them as exact strings.
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $item[2]

    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $item[3]

    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $item[4]

    $row++
}
